# Cleaned and Filtered Dataset and Added stuff to slides
#
# Adds three new "Title and Content" slides (indexes 2, 3, 4) after the
# existing title slide, each with a title + a bulleted content
# placeholder whose runs are colored dark red (C00000).

$p = $ppt.ActivePresentation

# Red accent color used throughout the new body text (RGB 192,0,0 ->
# OLE COLORREF packs as R + G*256 + B*65536 = 192).
$red = 192

# Builds up a content placeholder's text paragraph by paragraph so that
# each paragraph can get its own color / indent level without the
# "apply to whole shape" behavior that setting .IndentLevel / Font on a
# stale InsertAfter() range can trigger.
function Set-BodyParagraphs {
    param(
        $TextRange,
        [System.Collections.IEnumerable]$Lines
    )

    $i = 0
    foreach ($line in $Lines) {
        $text = $line[0]
        $level = $line[1]
        $i = $i + 1
        if ($i -eq 1) {
            $TextRange.Text = $text
        } else {
            $TextRange.InsertAfter("`r" + $text) | Out-Null
        }
        $para = $TextRange.Paragraphs($i, 1)
        $para.Font.Color.RGB = $red
        if ($level -gt 0) {
            $para.IndentLevel = $level + 1
        }
    }
}

# ---------------------------------------------------------------------
# Slide 2 - "Research Question/Objective"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Add(2, 2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Research Question/Objective"

$body2 = $s2.Shapes.Item(2).TextFrame.TextRange
Set-BodyParagraphs $body2 @(
    ,@("Why I Chose This Topic (Background)", 0)
    ,@("How it Addresses a Specific Problem or Question", 0)
)

# ---------------------------------------------------------------------
# Slide 3 - "Data Sources & Transformations"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Data Sources & Transformations"

$body3 = $s3.Shapes.Item(2).TextFrame.TextRange
Set-BodyParagraphs $body3 @(
    ,@("Where Data Comes From, Reliability, Limitations", 0)
    ,@("Summarize Cleaning and Transformation Steps in Non-Technical Way", 0)
)

# ---------------------------------------------------------------------
# Slide 4 - "Contextual Visualizations"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Add(4, 2)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Contextual Visualizations"

$tf4 = $s4.Shapes.Item(2).TextFrame
$body4 = $tf4.TextRange
Set-BodyParagraphs $body4 @(
    ,@("Bar Chart (Satellite Missions by Frequency Band)", 0)
    ,@("X = Mission Purpose", 1)
    ,@("Y = Number of Satellites", 1)
    ,@("Stacked Bar Chart  (Comparing Mission Types Across Multiple Bands)", 0)
    ,@("X = Mission Purpose", 1)
    ,@("Stacked = Frequency Bands", 1)
    ,@("Heatmap (Correlations Between Mission Type and Orbit Type)", 0)
    ,@("Rows: Mission Purpose", 1)
    ,@("Columns: Orbit Type ", 1)
    ,@("Color Intensity: Number of Satellites", 1)
)

$tf4.AutoSize = 2

Write-Output "done"
